$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New identifiers used by this handoff report generation.
# ---------------------------------------------------------------------------
$newFile1 = "158925e4-019d-4913-bfb5-2a29c1a7aa81"
$newFile2 = "ffff0948183c-f0f3-46cb-a263-07704fb4eb9a"
$newHash  = "a0f1ddf3db9d086a02c8310c1debe2f83fe237e5"

$status        = "Ready for handoff"
$overviewDate  = "2016-27-17 10:27:18"
$zhHandoffDate = "2016-03-17 10:27:15"
$deHandoffDate = "2016-03-17 10:27:18"
$emptyDate     = "0001-01-01 00:00:00"

# ===========================================================================
# Sheet "Overview"
# ===========================================================================
$ws = $wb.Worksheets.Item("Overview")
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "$newFile1.md"
$ws.Range("B2").Value = $status
$ws.Range("C2").Value = $status
$ws.Range("D2").Value = $overviewDate

$ws.Range("A3").Value = "$newFile2.md"
$ws.Range("B3").Value = $status
$ws.Range("C3").Value = $status
$ws.Range("D3").Value = $overviewDate

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f3ecc755f260fc6836a636eebb207ab60af5593d/e2e/$newFile1.md", [Type]::Missing, [Type]::Missing, "$newFile1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f3ecc755f260fc6836a636eebb207ab60af5593d/e2e/$newFile2.md", [Type]::Missing, [Type]::Missing, "$newFile2.md") | Out-Null

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Delete()

# Drop "Latest Target File" / "Latest Handback File" - no handback has happened yet.
$ws.Range("F2:G3").Clear()

$ws.Range("A2").Value = "$newFile1.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = $status
$ws.Range("D2").Value = "$newFile1.$newHash.zh-cn.xlf"
$ws.Range("E2").Value = $zhHandoffDate
$ws.Range("H2").Value = $emptyDate
$ws.Range("I2").Value = "Include"

$ws.Range("A3").Value = "$newFile2.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = $status
$ws.Range("D3").Value = "$newFile1.$newHash.zh-cn.xlf"
$ws.Range("E3").Value = $zhHandoffDate
$ws.Range("H3").Value = $emptyDate
$ws.Range("I3").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/58085a535928c143ce2a115d45dfbe9dfdb6925b/e2e/$newFile1.md", [Type]::Missing, [Type]::Missing, "$newFile1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/58085a535928c143ce2a115d45dfbe9dfdb6925b/e2e/$newFile1.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/265b2c67abacd95070db083ae99917fa40db4a21/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newFile1.$newHash.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "$newFile1.$newHash.zh-cn.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/58085a535928c143ce2a115d45dfbe9dfdb6925b/e2e/$newFile2.md", [Type]::Missing, [Type]::Missing, "$newFile2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/58085a535928c143ce2a115d45dfbe9dfdb6925b/e2e/$newFile2.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/265b2c67abacd95070db083ae99917fa40db4a21/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newFile1.$newHash.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "$newFile1.$newHash.zh-cn.xlf") | Out-Null

# ===========================================================================
# Sheet "de-de"
# ===========================================================================
$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Delete()

# Drop "Latest Target File" / "Latest Handback File" - no handback has happened yet.
$ws.Range("F2:G3").Clear()

$ws.Range("A2").Value = "$newFile1.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = $status
$ws.Range("D2").Value = "$newFile1.$newHash.de-de.xlf"
$ws.Range("E2").Value = $deHandoffDate
$ws.Range("H2").Value = $emptyDate
$ws.Range("I2").Value = "Include"

$ws.Range("A3").Value = "$newFile2.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = $status
$ws.Range("D3").Value = "$newFile1.$newHash.de-de.xlf"
$ws.Range("E3").Value = $deHandoffDate
$ws.Range("H3").Value = $emptyDate
$ws.Range("I3").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/363ae2ba11fb733b4529e6a9e06b536dc13ea9fa/e2e/$newFile1.md", [Type]::Missing, [Type]::Missing, "$newFile1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/363ae2ba11fb733b4529e6a9e06b536dc13ea9fa/e2e/$newFile1.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/811854fef82c7b39bc1827d1cd544e7b53af3d40/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newFile1.$newHash.de-de.xlf", [Type]::Missing, [Type]::Missing, "$newFile1.$newHash.de-de.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/363ae2ba11fb733b4529e6a9e06b536dc13ea9fa/e2e/$newFile2.md", [Type]::Missing, [Type]::Missing, "$newFile2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/363ae2ba11fb733b4529e6a9e06b536dc13ea9fa/e2e/$newFile2.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/811854fef82c7b39bc1827d1cd544e7b53af3d40/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newFile1.$newHash.de-de.xlf", [Type]::Missing, [Type]::Missing, "$newFile1.$newHash.de-de.xlf") | Out-Null

Write-Host "Report regenerated for handoff."
